$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange
$para1 = $tr.Paragraphs(1,1)
Write-Host "ParagraphFormat: $($para1.ParagraphFormat)"
try {
  Write-Host "Indent: $($para1.IndentLevel)"
} catch { Write-Host "err1 $_" }
try {
  Write-Host "LeftIndent: $($para1.ParagraphFormat.LeftIndent)"
} catch { Write-Host "err2 $_" }
